$d = $word.ActiveDocument

# The document starts as a single paragraph made of 4 runs whose combined
# text reads "First demo 3" (runs: "First ", "demo", " ", "3"), followed by
# the _GoBack bookmark. We need three paragraphs reading "Use case 1",
# "Use case 2" and "Use case 3", with the bookmark still trailing the last
# run of the final paragraph.

# Step 1: collapse all the runs' text into the final combined wording while
# keeping everything (including the bookmark) inside a single paragraph.
$d.Content.Find.Execute("First demo 3", $true, $false, $false, $false, $false, $true, 1, $false, "Use case 1Use case 2Use case 3", 2)

# Step 2: split that single paragraph into three, by inserting paragraph
# breaks at the boundaries between "Use case 1" | "Use case 2" | "Use case 3".
# Insert from the end backwards so earlier character offsets stay valid.
$d.Range(20, 20).InsertParagraphBefore()
$d.Range(10, 10).InsertParagraphBefore()
